$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.618.27"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.530.92"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "317.25"
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("D6").Value = "94.86"
$ws.Range("E6").Value = "  -6.18%  "
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  -4.18%  "
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "2.919.21"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "15.50"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.526.22"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "0.848"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").Value = "42.645.08"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").Value = "6.56"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").Value = "70.11"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").Value = "251.14"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "26.48"
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  +2.42%  "
$ws.Range("D29").Value = "39.14"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").Value = "10.15"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "6.02"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").Value = "155.83"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").Value = "19.31"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("D34").Value = "2.12"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "3.28"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "0.0783"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "0.110"
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").Value = "23.74"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  +10.11%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "3.80"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "0.0301"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  -6.02%  "
$ws.Range("D46").Value = "2.010.42"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").Value = "85.51"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "8.80"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").Value = "2.774.37"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "74.38"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").Value = "102.61"
$ws.Range("E51").Value = "  -0.87%  "